# Generate Report for Handoff
# Update the localization-status workbook to reflect that the
# c7a2c94e-... file is now "Ready for handoff" and refresh the
# corresponding handoff timestamps on the language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the c7a2c94e-...md file; mark both
# language statuses as "Ready for handoff"
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the c7a2c94e-...md file
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D2").Value = "2016-02-18 10:32:08"
$zhcn.Range("D3").Value = "2016-02-18 10:32:08"

# de-de sheet: row 3 is the c7a2c94e-...md file
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D2").Value = "2016-02-18 10:32:19"
$dede.Range("D3").Value = "2016-02-18 10:32:19"
